# Update Mitelmandb version string on the "biomarkers" sheet (row 3, column E:
# source_version for "Mitelman Database" / mitelmandb), from v20250710 to v20250815.
#
# Also reflect that, after making this edit, the "biomarkers" sheet/tab and cell E3
# are the active selection when the workbook was saved.

$wb = $excel.ActiveWorkbook

$wsCompounds  = $wb.Worksheets.Item("compounds")
$wsBiomarkers = $wb.Worksheets.Item("biomarkers")

# Update the Mitelman Database version value.
$wsBiomarkers.Range("E3").Value = "v20250815"

# Make "biomarkers" the active sheet, with E3 selected, to match the saved view state.
$wsBiomarkers.Activate()
$wsBiomarkers.Range("E3").Select()
